# Scheduled runner update: refresh Leve profit market-price computations across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 28788
$ws.Range("I6").Value = 28788
$ws.Range("K6").Value = 86364
$ws.Range("M6").Value = -86252

$ws.Range("H32").Value = 4618.8887
$ws.Range("I32").Value = 4341.4287
$ws.Range("J32").Value = 4795.4546
$ws.Range("K32").Value = 4341.4287
$ws.Range("L32").Value = 4795.4546
$ws.Range("M32").Value = -4015.4287
$ws.Range("N32").Value = -5447.4546

$ws.Range("H33").Value = 2160.2727
$ws.Range("I33").Value = 1846.1177
$ws.Range("K33").Value = 1846.1177
$ws.Range("M33").Value = -1617.1177

$ws.Range("H40").Value = 3602.9092
$ws.Range("I40").Value = 3309.2307
$ws.Range("K40").Value = 3309.2307
$ws.Range("M40").Value = -3134.2307

$ws.Range("H53").Value = 1212.2727
$ws.Range("I53").Value = 408.75
$ws.Range("K53").Value = 408.75
$ws.Range("M53").Value = 228.25

$ws.Range("H112").Value = 2110.9443
$ws.Range("J112").Value = 2193.8667
$ws.Range("L112").Value = 6581.6001
$ws.Range("N112").Value = -8797.6001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 43099.8
$ws.Range("J44").Value = 43099.8
$ws.Range("L44").Value = 43099.8
$ws.Range("N44").Value = -44075.8

$ws.Range("H45").Value = 31252280
$ws.Range("I45").Value = 41668548
$ws.Range("K45").Value = 41668548
$ws.Range("M45").Value = -41668171

$ws.Range("H55").Value = 2999
$ws.Range("I55").Value = 2999
$ws.Range("K55").Value = 2999
$ws.Range("M55").Value = -2684

$ws.Range("H74").Value = 11911014
$ws.Range("I74").Value = 14707370
$ws.Range("K74").Value = 14707370
$ws.Range("M74").Value = -14706496

$ws.Range("H77").Value = 11911014
$ws.Range("I77").Value = 14707370
$ws.Range("K77").Value = 73536850
$ws.Range("M77").Value = -73532482

$ws.Range("H96").Value = 34997.5
$ws.Range("J96").Value = 34997.5
$ws.Range("L96").Value = 34997.5
$ws.Range("N96").Value = -40489.5

$ws.Range("H133").Value = 66998.625
$ws.Range("J133").Value = 66570.28999999999
$ws.Range("L133").Value = 66570.28999999999
$ws.Range("N133").Value = -71630.28999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 3840.6667
$ws.Range("I22").Value = 2208.8
$ws.Range("K22").Value = 2208.8
$ws.Range("M22").Value = -2035.8

$ws.Range("H54").Value = 16759.6
$ws.Range("J54").Value = 6899.5
$ws.Range("L54").Value = 6899.5
$ws.Range("N54").Value = -7867.5

$ws.Range("H105").Value = 2703.1304
$ws.Range("I105").Value = 1687.8
$ws.Range("J105").Value = 2985.1667
$ws.Range("K105").Value = 1687.8
$ws.Range("L105").Value = 2985.1667
$ws.Range("M105").Value = 59.20000000000005
$ws.Range("N105").Value = -6479.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value = 56817.8
$ws.Range("J53").Value = 56817.8
$ws.Range("L53").Value = 56817.8
$ws.Range("N53").Value = -58031.8

$ws.Range("H107").Value = 14454.667
$ws.Range("I107").Value = 17682
$ws.Range("J107").Value = 8000
$ws.Range("K107").Value = 17682
$ws.Range("L107").Value = 8000
$ws.Range("M107").Value = -15762
$ws.Range("N107").Value = -11840

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 196.07692
$ws.Range("I8").Value = 196.07692
$ws.Range("K8").Value = 588.23076
$ws.Range("M8").Value = -449.23076

$ws.Range("H14").Value = 848.7143
$ws.Range("I14").Value = 848.7143
$ws.Range("K14").Value = 2546.1429
$ws.Range("M14").Value = -2373.1429

$ws.Range("H37").Value = 64825
$ws.Range("J37").Value = 64825
$ws.Range("L37").Value = 194475
$ws.Range("N37").Value = -194699

$ws.Range("H82").Value = 9999.833000000001
$ws.Range("I82").Value = 9999
$ws.Range("K82").Value = 29997
$ws.Range("M82").Value = -29591

$ws.Range("H85").Value = 9999.833000000001
$ws.Range("I85").Value = 9999
$ws.Range("K85").Value = 29997
$ws.Range("M85").Value = -28593

$ws.Range("H104").Value = 2449.75
$ws.Range("I104").Value = 899.5
$ws.Range("K104").Value = 2698.5
$ws.Range("M104").Value = -77.5

$ws.Range("H133").Value = 5000
$ws.Range("I133").Value = 5000
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 15000
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -9940
$ws.Range("N133").ClearContents()

$ws.Range("H140").Value = 190473.56
$ws.Range("I140").Value = 190473.56
$ws.Range("K140").Value = 571420.6799999999
$ws.Range("M140").Value = -566240.6799999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 315000
$ws.Range("I19").Value = 425000
$ws.Range("K19").Value = 425000
$ws.Range("M19").Value = -424712

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()

$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()

$ws.Range("H68").Value = 1960.4546
$ws.Range("I68").Value = 1674
$ws.Range("K68").Value = 1674
$ws.Range("M68").Value = -925

$ws.Range("H71").Value = 1960.4546
$ws.Range("I71").Value = 1674
$ws.Range("K71").Value = 8370
$ws.Range("M71").Value = -4626

$ws.Range("H109").Value = 47500
$ws.Range("J109").Value = 47500
$ws.Range("L109").Value = 47500
$ws.Range("N109").Value = -50274

$ws.Range("H132").Value = 444472.62
$ws.Range("I132").Value = 97012.45
$ws.Range("K132").Value = 291037.35
$ws.Range("M132").Value = -288507.35

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 83343.336
$ws.Range("J24").Value = 83343.336
$ws.Range("L24").Value = 83343.336
$ws.Range("N24").Value = -83803.336

$ws.Range("H82").Value = 32996
$ws.Range("J82").Value = 34999
$ws.Range("L82").Value = 34999
$ws.Range("N82").Value = -35765

$ws.Range("H85").Value = 32996
$ws.Range("J85").Value = 34999
$ws.Range("L85").Value = 34999
$ws.Range("N85").Value = -37651

$ws.Range("H88").Value = 18000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 18000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 18000
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -18812

$ws.Range("H91").Value = 18000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 18000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 18000
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -20808

$ws.Range("H99").Value = 111000
$ws.Range("J99").Value = 111000
$ws.Range("L99").Value = 111000
$ws.Range("N99").Value = -116990

$ws.Range("H132").Value = 437090.6
$ws.Range("I132").Value = 2412.682
$ws.Range("K132").Value = 7238.045999999999
$ws.Range("M132").Value = -4708.045999999999

$ws.Range("H136").Value = 5212.2856
$ws.Range("J136").Value = 4173.7144
$ws.Range("L136").Value = 12521.1432
$ws.Range("N136").Value = -17621.1432
